$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# New row 11 values (mirrors the layout of the existing TestCase rows)
$ws.Range("A11").Value = "TestCase10_validateFecthTableData"
$ws.Range("B11").Value = "Validate different column values of different Candidates in the web table."
$ws.Range("C11").Value = "kw_fetchtabledata"
$ws.Range("D11").Value = "Y"
$ws.Range("J11").Value = "Table"
$ws.Range("K11").Value = "Table Data Download"
$ws.Range("L11").Value = "Airi Satou#Bradley Greer#Brenden Wagner#Colleen Hurst"
$ws.Range("N11").Value = "Position#Age#Salary"

# Match the formatting (borders/fills/wrap) of the row above it
$ws.Range("A10:N10").Copy()
$ws.Range("A11:N11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row height matches the wrapped-text row above
$ws.Rows.Item(11).RowHeight = 60

# Update the view's scroll position / active selection
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("M11").Select()
